$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The nightly export from Scholarship Universe wrote TotalAmountAvailable
# (column B) out as text, and the rows within each scholarship group were
# left in the order the extract happened to return them. Put the amounts
# back as real numbers and re-sort each group by ApplicantRanking (column D)
# ascending, the way the new clean-up stored procedure will produce them.

$groups = @(
    @(2, 24, 12000),
    @(25, 33, 4500),
    @(34, 44, 4000),
    @(45, 48, 4000),
    @(49, 50, 1500),
    @(51, 58, 1000)
)

foreach ($g in $groups) {
    $firstRow = $g[0]
    $lastRow = $g[1]
    $amount = $g[2]

    # Re-sort this scholarship's rows by ApplicantRanking ascending.
    $dataRange = $ws.Range("A$firstRow`:D$lastRow")
    $keyRange = $ws.Range("D$firstRow`:D$lastRow")
    $dataRange.Sort($keyRange)

    # Replace the text TotalAmountAvailable with the real number.
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 2).Value = $amount
    }
}

$ws.Range("A1:D58").Select() | Out-Null
